# RA-5 False positive findings are documented
#
# The CivicActions Security Office paragraph (FirstParagraph style) gets a
# new sentence inserted between "...security assessments." and
# "Vulnerabilities found and deemed legitimate...":
#
#   ... security assessments. \u201cFalse positive\u201d findings are
#   documented and may be tailored out. Vulnerabilities found ...

$d = $word.ActiveDocument

$old = "The CivicActions Security Office reviews all vulnerabilities identified from automated scans and security assessments. Vulnerabilities found and deemed legitimate are assigned an impact rating and response time thought creation of an issue or ticket. The CivicActions Operations staff reviews current scans and compare with older scans to identify trends and to verify previous vulnerabilities have been mitigated."

$new = "The CivicActions Security Office reviews all vulnerabilities identified from automated scans and security assessments. " + [char]0x201C + "False positive" + [char]0x201D + " findings are documented and may be tailored out. Vulnerabilities found and deemed legitimate are assigned an impact rating and response time thought creation of an issue or ticket. The CivicActions Operations staff reviews current scans and compare with older scans to identify trends and to verify previous vulnerabilities have been mitigated."

$r = $d.Content
$found = $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
Write-Host "Found/replaced: $found"
